$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp string (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 10 de Octubre de 2020 a las 14:48"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 7895738
$ws.Range("C4").Value = 1260
$ws.Range("D4").Value = 5065384
$ws.Range("E4").Value = 2611669
$ws.Range("G4").Value = 37
$ws.Range("H4").Value = 218685

# --- Row 25: Alemania ---
$ws.Range("D25").Value = 273500
$ws.Range("E25").Value = 37291

# --- Row 30: Paises Bajos ---
$ws.Range("B30").Value = 168280
$ws.Range("C30").Value = 6499
$ws.Range("G30").Value = 23
$ws.Range("H30").Value = 6567

# --- Rows 79/80: Bosnia y Herzegovina / El Salvador swap places & update data ---
# Row 79 keeps the higher total (now Bosnia y Herzegovina, freshly updated figures).
# Row 80 becomes El Salvador, carrying the figures that used to sit in row 79.
$ws.Range("A79").Value = "Bosnia y Herzegovina"
$ws.Range("B79").Value = 30345
$ws.Range("C79").Value = 428
$ws.Range("D79").Value = 23370
$ws.Range("E79").Value = 6048
$ws.Range("G79").Value = 1
$ws.Range("H79").Value = 927

$ws.Range("A80").Value = "El Salvador"
$ws.Range("B80").Value = 29951
$ws.Range("C80").Value = 0
$ws.Range("D80").Value = 24995
$ws.Range("E80").Value = 4069
$ws.Range("G80").Value = 6
$ws.Range("H80").Value = 887

# --- Row 96: Senegal ---
$ws.Range("B96").Value = 15244
$ws.Range("C96").Value = 31
$ws.Range("D96").Value = 13198
$ws.Range("E96").Value = 1732
$ws.Range("G96").Value = 1
$ws.Range("H96").Value = 314

# --- Row 127: Hong Kong ---
$ws.Range("D127").Value = 4914
$ws.Range("E127").Value = 157

# --- Row 136: Sri Lanka ---
$ws.Range("D136").Value = 3306
$ws.Range("E136").Value = 1204

# --- Row 168: Vietnam ---
$ws.Range("B168").Value = 1107
$ws.Range("C168").Value = 2
$ws.Range("E168").Value = 48
